# LOT14IMEIGUARD - envio de imagem PM
# Adds three new form-response rows (11-13) to the "Respostas ao formulário 1"
# sheet, matching the formatting of the existing response rows, plus the
# associated page setup / selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Formatting: clone the look of existing response rows ----------------
# Row 11 mirrors row 10's shape (13.5pt row, bottom border).
$ws.Range("A10:J10").Copy()
$ws.Range("A11:J11").PasteSpecial(-4122)  # xlPasteFormats

# Rows 12-13 mirror row 9's shape (26.25pt row, bottom border).
$ws.Range("A9:J9").Copy()
$ws.Range("A12:J12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A13:J13").PasteSpecial(-4122)  # xlPasteFormats

$ws.Rows.Item(11).RowHeight = 13.5
$ws.Rows.Item(12).RowHeight = 26.25
$ws.Rows.Item(13).RowHeight = 26.25

# Matricula/CPF columns (D/E) need the right-aligned numeric look on some
# rows and the plain "general" look on others, matching how each response
# was originally keyed in.
$ws.Range("D11").HorizontalAlignment = -4152  # xlRight
$ws.Range("E11").HorizontalAlignment = -4152  # xlRight
$ws.Range("D12").HorizontalAlignment = -4152  # xlRight
$ws.Range("E12").HorizontalAlignment = -4152  # xlRight
$ws.Range("D13").HorizontalAlignment = -4152  # xlRight
$ws.Range("E13").HorizontalAlignment = 1      # xlGeneral

# --- 2. Cell values -----------------------------------------------------
# Text values are entered in the same order the three new form responses
# were originally captured so repeated values (e.g. the same e-mail typed
# into two columns) are reused rather than duplicated.
$ws.Range("B11").Value = "direito.ariclessilva@gmail.com"
$ws.Range("F11").Value = "MPC"
$ws.Range("B12").Value = "jacksonwla@gmail.com"
$ws.Range("C12").Value = "JACKSON WENDELL LOPES DE ALMEIDA"
$ws.Range("G12").Value = "PM"
$ws.Range("H12").Value = "ROTAM"
$ws.Range("B13").Value = "jefferson.rodrisouza@gmail.com"
$ws.Range("I13").Value = "jefferson.souza@policiacivil.pa.gov.br"
$ws.Range("C11").Value = "ARICLES DE SOUSA SILVA"
$ws.Range("C13").Value = "JEFFERSON RODRIGUES SOUZA"
$ws.Range("H11").Value = "SUPERINTENDENCIA CASTANHAL, 3 RISP"
$ws.Range("H13").Value = "DELEGACIA DE BENFICA"
$ws.Range("F12").Value = "SGT"
$ws.Range("I11").Value = "direito.ariclessilva@gmail.com"
$ws.Range("I12").Value = "jacksonwla@gmail.com"
$ws.Range("G11").Value = "PC"
$ws.Range("J11").Value = "CONFIRMADO"
$ws.Range("G13").Value = "PC"
$ws.Range("F13").Value = "IPC"
$ws.Range("J12").Value = "CONFIRMADO"
$ws.Range("J13").Value = "CONFIRMADO"

# Timestamps (Carimbo de data/hora)
$ws.Range("A11").Value = 45485.94604166667
$ws.Range("A12").Value = 45488.699837962966
$ws.Range("A13").Value = 45488.726354166669

# Matricula funcional / CPF
$ws.Range("D11").Value = 5453003
$ws.Range("E11").Value = 17116414215
$ws.Range("D12").Value = 541930761
$ws.Range("E12").Value = 72273275234
$ws.Range("D13").Value = 4219631
$ws.Range("E13").Value = 2041287225

# --- 3. Selection ---------------------------------------------------------
$ws.Range("A11:J13").Select()

# --- 4. Page setup ---------------------------------------------------------
$ps = $ws.PageSetup
$ps.PaperSize = 9       # xlPaperA4
$ps.Orientation = 1     # xlPortrait
